# Applies the "verbose_name" column insertion + True/False text conversion
# to the Customer.xlsx model workbook, plus the Meta-class shared-string
# edit and the small view/selection/row-height tweaks.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "model" sheet: insert a new column D ("verbose_name") and populate it.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("model")

# Inserting a column at D pushes the old D:W block to E:X, adjusting the
# dimension, shared cols and every cell reference automatically.
$ws1.Columns("D").Insert()

# Header for the freshly inserted column.
$ws1.Range("D1").Value2 = "verbose_name"

# Per-row verbose_name literal values (quoted python/django style strings).
# D2 is typed with a leading apostrophe in the source workbook, so it picks
# up the text "quote prefix" cell style; the rest do not.
$ws1.Range("D2").Value2 = "'" + '"Company Name"'
$ws1.Range("D3").Value2 = '"Active"'
$ws1.Range("D4").Value2 = '"AFM"'
$ws1.Range("D5").Value2 = '"First Name"'
$ws1.Range("D6").Value2 = '"Last Name"'
$ws1.Range("D7").Value2 = '"Email"'
$ws1.Range("D8").Value2 = '"Phone"'
$ws1.Range("D9").Value2 = '"Address"'
$ws1.Range("D10").Value2 = '"Created at"'
$ws1.Range("D11").Value2 = '"Updated at"'

# The AFM row used to carry a numeric 9 in its "max_length" column (now L);
# that value is gone in the edited workbook.
$ws1.Range("L4").ClearContents()

# Every boolean TRUE/FALSE flag cell became literal text "True"/"False"
# (typed with a leading apostrophe, hence the quote-prefix style).
$trueCells = @("T2","U2","V2","X2","T3","V3","X3","T4","V4","X4","T5","V5","X5","T6","V6","X6","T7","V7","X7","T8","V8","X8","T9","V9","X9","T10","V10","T11","V11","N4","O4","N8","O8","N9","O9","M7","Q10","R11")
foreach ($cellRef in $trueCells) {
    $ws1.Range($cellRef).Value2 = "'True"
}
$ws1.Range("P3").Value2 = "'False"

# X10 / X11 end up as blank cells that still carry the quote-prefix style
# (set then clear, so the style sticks but the content disappears).
$ws1.Range("X10").Value2 = "'True"
$ws1.Range("X10").ClearContents()
$ws1.Range("X11").Value2 = "'True"
$ws1.Range("X11").ClearContents()

# ---------------------------------------------------------------------
# 2) "model_functions" sheet: update the Meta class snippet + view state.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("model_functions")
$ws2.Range("B3").Value2 = "    class Meta:`n        verbose_name = 'Customer'`n        verbose_name_plural = 'Customers'"
$ws2.Rows(3).RowHeight = 51
$ws2.Range("B10").Select()

# ---------------------------------------------------------------------
# 3) Restore "model" as the active sheet/selection last, so it keeps the
#    tab-selected view state and picks up its new cursor position.
# ---------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("L4").Select()

Write-Output "done"
